$d = $word.ActiveDocument

# The document currently has several short paragraphs after the first one
# ("Hello"+"dsldksld", then a handful of filler paragraphs). The edit keeps
# only the first paragraph (retitled) and the final (empty) paragraph that
# precedes the section properties.

$firstPara = $d.Paragraphs.Item(1)
$lastPara  = $d.Paragraphs.Item($d.Paragraphs.Count)

# Delete everything between the end of the first paragraph and the start of
# the last paragraph (i.e. all the in-between paragraphs, including their
# paragraph marks).
$middle = $d.Range($firstPara.Range.End, $lastPara.Range.Start)
$middle.Delete()

# Replace the first paragraph's text (currently "Hello" + "dsldksld") with
# the new wording, leaving the paragraph mark untouched.
$firstPara = $d.Paragraphs.Item(1)
$textRange = $d.Range($firstPara.Range.Start, $firstPara.Range.End - 1)
$textRange.Text = "Nouvelle version"
